$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: a second "location" entry (en-dash variant) is added next to the
# existing one, in a new column E2 - reuse the same look/format that the
# neighbouring "language" cell (G2) already had.
$null = $ws.Range("G2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E2").Value = "../../QML/OggettiEditDash/AddNewForm.qml – 21"

# Row 1 got a touch shorter.
$ws.Rows.Item(1).RowHeight = 13.8

# Selection moved.
$null = $ws.Range("E8").Select()
